$d = $word.ActiveDocument

# --- Helper: replace the text of a whole paragraph (excluding the
# trailing paragraph mark) with a new string. Using a single Range.Text
# assignment per paragraph/block keeps the edit atomic and predictable. ---

# 1) Title
$p1 = $d.Paragraphs(1)
$r1 = $d.Range($p1.Range.Start, $p1.Range.End - 1)
$r1.Text = "Politics: Navigating the Maze of Power and Influence"

# 2) Author name line: "Isabella J. Hutchinson" -> "Emma Watson"
$p2 = $d.Paragraphs(2)
$r2 = $d.Range($p2.Range.Start, $p2.Range.End - 1)
$r2.Text = "Emma Watson"

# 3) Author email line: "isabella.hutchinson@berkeley.edu" -> "emma.watson87@schoolmail.net"
$p3 = $d.Paragraphs(3)
$r3 = $d.Range($p3.Range.Start, $p3.Range.End - 1)
$r3.Text = "emma.watson87@schoolmail.net"

# 4) Paragraph 4 is an already-empty paragraph; leave untouched.

# 5) Main body paragraph (three sub-sections separated by double line breaks)
$p5 = $d.Paragraphs(5)
$r5 = $d.Range($p5.Range.Start, $p5.Range.End - 1)

$br2 = [string][char]11 + [char]11

$block1 = "Politics, a multifaceted and dynamic realm of human interaction, permeates every aspect of our lives." + `
  " It shapes the laws, policies, and decisions that govern societies, impacting individuals, communities, and nations alike." + `
  " As citizens of a democratic society, it is imperative for us to understand the intricacies of politics and the role we play in shaping its course."

$block2 = "Politics is often perceived as a complex web of power dynamics, negotiations, and compromises." + `
  " It involves the interactions among various stakeholders, including elected officials, political parties, interest groups, and the general public." + `
  " Understanding the different branches of government, their functions, and how they interact is crucial for comprehending the political landscape." + `
  " Political ideologies, such as liberalism, conservatism, and socialism, influence the policy positions and actions of political actors, and it is essential to grasp these ideologies and their implications."

$block3 = "Beyond the formal institutions and processes, politics also encompasses the informal dynamics of influence and persuasion." + `
  " Lobbying, public relations, and grassroots movements play a significant role in shaping political outcomes." + `
  " The media plays a vital role in informing and shaping public opinion, and understanding the relationship between politics and the media is crucial for informed citizenship." + `
  " The influence of money in politics, campaign finance regulations, and the role of special interest groups are important aspects to consider in analyzing the political landscape."

$r5.Text = $block1 + $br2 + $block2 + $br2 + $block3

# 6) "Summary" heading paragraph stays the same.

# 7) Summary body paragraph
$p7 = $d.Paragraphs(7)
$r7 = $d.Range($p7.Range.Start, $p7.Range.End - 1)
$r7.Text = "Politics, a complex and ever-evolving field, involves the interactions among various stakeholders, encompassing formal institutions, informal dynamics, and the influence of the media." + `
  " It requires an understanding of political ideologies, government structures, and the role of interest groups." + `
  " As citizens, it is essential for us to engage with politics, be informed about current issues, and participate in the democratic process to create a society that reflects our values and aspirations."

# 8) Append a new empty paragraph at the very end of the document body.
$endPos = $d.Content.End
$rEnd = $d.Range($endPos, $endPos)
$rEnd.InsertParagraphAfter()
